$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"="1.02"; "C"="1.10747583238813"; "D"="1.109065147348365"; "E"="1.106216010620477"; "F"="1.114365885490917"; "I"="1.031498234603447"; "J"="1.112224972927905"; "K"="1.111665833271221"; "L"="1.108823755236945"; "M"="1.11695354222363"; "N"="1.11380445996345" }
    3 = @{ "B"="1.02"; "C"="1.110857674634002"; "D"="1.112198861315344"; "E"="1.109247021146085"; "F"="1.11758243941636"; "I"="1.031820444163201"; "J"="1.115267425907758"; "K"="1.114617726851672"; "L"="1.111672646204427"; "M"="1.119989075358422"; "N"="1.116851233575509" }
    4 = @{ "B"="1.02"; "C"="1.113019560727219"; "D"="1.114200834055735"; "E"="1.111183350419882"; "F"="1.119633968715179"; "I"="1.032018224295869"; "J"="1.117210015193638"; "K"="1.116501895288911"; "L"="1.113490958547457"; "M"="1.121923336904248"; "N"="1.118796581560993" }
    5 = @{ "B"="1.02"; "C"="1.11392226319464"; "D"="1.115036452283791"; "E"="1.111991561452339"; "F"="1.120489463005397"; "I"="1.032098830286184"; "J"="1.118020587778702"; "K"="1.1172879467574"; "L"="1.114249509059834"; "M"="1.122729494183184"; "N"="1.119608305252102" }
    6 = @{ "B"="1.02"; "C"="1.114073474660707"; "D"="1.115176408118562"; "E"="1.112126926452435"; "F"="1.12063270028322"; "I"="1.03211221604978"; "J"="1.118156333556848"; "K"="1.11741957752257"; "L"="1.1143765329965"; "M"="1.122864445200891"; "N"="1.119744243804823" }
    7 = @{ "B"="1.02"; "C"="1.113031646650869"; "D"="1.114212023046912"; "E"="1.111194172457661"; "F"="1.119645427039146"; "I"="1.032019311310106"; "J"="1.117220869843297"; "K"="1.116512422137381"; "L"="1.113501117206632"; "M"="1.121934136139504"; "N"="1.118807451625499" }
    8 = @{ "B"="1.02"; "C"="1.108624307087935"; "D"="1.110129626587055"; "E"="1.107245609314728"; "F"="1.115459205677604"; "I"="1.03160935919392"; "J"="1.113258677828044"; "K"="1.112668894442869"; "L"="1.109791838815323"; "M"="1.117985710356761"; "N"="1.114839632843115" }
    9 = @{ "B"="1.02"; "C"="1.100648210693009"; "D"="1.102731522365761"; "E"="1.100089803111166"; "F"="1.107846620736463"; "I"="1.030803809455026"; "J"="1.106069917850919"; "K"="1.105690750142487"; "L"="1.103056571024624"; "M"="1.110791359663047"; "N"="1.10764066400226" }
    10 = @{ "B"="1.02"; "C"="1.095178559385279"; "D"="1.097651526156541"; "E"="1.095176042395697"; "F"="1.102601637729324"; "I"="1.03020917183332"; "J"="1.10112791892891"; "K"="1.100890406594694"; "L"="1.098422732192612"; "M"="1.105825017059328"; "N"="1.102691646874927" }
    11 = @{ "B"="1.02"; "C"="1.092771391597898"; "D"="1.095414271222653"; "E"="1.093011968644956"; "F"="1.100287529848434"; "I"="1.029937613795872"; "J"="1.098950068944764"; "K"="1.098774242956535"; "L"="1.096379835774215"; "M"="1.103631601207713"; "N"="1.100510704094" }
    12 = @{ "B"="1.02"; "C"="1.091871200329811"; "D"="1.094577386249764"; "E"="1.092202454254562"; "F"="1.099421267338795"; "I"="1.029834592565806"; "J"="1.098135199134767"; "K"="1.097982343268249"; "L"="1.095615333777759"; "M"="1.102810182102769"; "N"="1.099694677075406" }
    13 = @{ "B"="1.02"; "C"="1.092064572659245"; "D"="1.094757170188492"; "E"="1.092376358485656"; "F"="1.099607390757046"; "I"="1.029856789004804"; "J"="1.098310262999106"; "K"="1.098152477303664"; "L"="1.095779582514214"; "M"="1.102986685751403"; "N"="1.099869989550513" }
    14 = @{ "B"="1.02"; "C"="1.092697106578651"; "D"="1.095345215097212"; "E"="1.092945171103858"; "F"="1.100216062423763"; "I"="1.029929142165031"; "J"="1.098882833618013"; "K"="1.098708905155282"; "L"="1.096316758898625"; "M"="1.103563840244992"; "N"="1.100443373285374" }
    15 = @{ "B"="1.02"; "C"="1.093086021057125"; "D"="1.095706744454911"; "E"="1.093294875933609"; "F"="1.100590189952561"; "I"="1.029973434943167"; "J"="1.09923482212441"; "K"="1.099050955275323"; "L"="1.096646972024865"; "M"="1.10391855094653"; "N"="1.100795861655823" }
    16 = @{ "B"="1.02"; "C"="1.095337478332475"; "D"="1.09779919482401"; "E"="1.095318880188491"; "F"="1.102754291345926"; "I"="1.030226894512538"; "J"="1.101271637899493"; "K"="1.101030039386662"; "L"="1.098557527734568"; "M"="1.105969662168249"; "N"="1.102835569942932" }
    17 = @{ "B"="1.02"; "C"="1.096739213818591"; "D"="1.099101517930718"; "E"="1.096578595362848"; "F"="1.104100097946788"; "I"="1.030382088538029"; "J"="1.102538970851501"; "K"="1.1022612553789"; "L"="1.099746075821565"; "M"="1.10724460840437"; "N"="1.104104702653097" }
    18 = @{ "B"="1.02"; "C"="1.097553097994781"; "D"="1.099857530328566"; "E"="1.097309870541209"; "F"="1.10488095255173"; "I"="1.030471255125722"; "J"="1.103274540298478"; "K"="1.102975791725912"; "L"="1.100435837039599"; "M"="1.107984134611951"; "N"="1.104841316693061" }
    19 = @{ "B"="1.02"; "C"="1.097829986744222"; "D"="1.10011470516991"; "E"="1.097558629958568"; "F"="1.105146510098162"; "I"="1.030501429825703"; "J"="1.103524738833968"; "K"="1.10321882402896"; "L"="1.100670440503527"; "M"="1.108235600581974"; "N"="1.105091870539168" }
    20 = @{ "B"="1.02"; "C"="1.096589207776606"; "D"="1.098962166025406"; "E"="1.096443802909609"; "F"="1.103956134686381"; "I"="1.030365578155125"; "J"="1.102403376539563"; "K"="1.102129532596097"; "L"="1.099618919230093"; "M"="1.107108247402946"; "N"="1.103968915781683" }
    21 = @{ "B"="1.02"; "C"="1.092511010404611"; "D"="1.095172214379103"; "E"="1.092777828568542"; "F"="1.100037010725214"; "I"="1.029907895696273"; "J"="1.098714391070324"; "K"="1.098545214636258"; "L"="1.096158732345445"; "M"="1.103394069150196"; "N"="1.100274691529946" }
    22 = @{ "B"="1.02"; "C"="1.089911687483537"; "D"="1.092755247742033"; "E"="1.090439902122848"; "F"="1.097534011248427"; "I"="1.029607658605537"; "J"="1.096360614810129"; "K"="1.096257579627615"; "L"="1.093950206019811"; "M"="1.101020006214189"; "N"="1.097917572637479" }
    23 = @{ "B"="1.02"; "C"="1.091293056228168"; "D"="1.094039834398723"; "E"="1.091682482078913"; "F"="1.098864668352542"; "I"="1.029768015552732"; "J"="1.097611729983918"; "K"="1.097473598925341"; "L"="1.095124185070369"; "M"="1.102282301657444"; "N"="1.09917046453832" }
    24 = @{ "B"="1.02"; "C"="1.096657000543414"; "D"="1.099025144295798"; "E"="1.096504720604177"; "F"="1.104021198238322"; "I"="1.030373042672571"; "J"="1.102464657019756"; "K"="1.102189063596287"; "L"="1.099676386618309"; "M"="1.107169875811678"; "N"="1.104030283287188" }
    25 = @{ "B"="1.02"; "C"="1.102736205480099"; "D"="1.104669382301452"; "E"="1.101964221290468"; "F"="1.109843744704099"; "I"="1.030803809455026"; "J"="1.107953940420296"; "K"="1.107520121260354"; "L"="1.104822370404422"; "M"="1.112680418627366"; "N"="1.109527362099786" }
}

foreach ($rowKey in $data.Keys) {
    $rowData = $data[$rowKey]
    foreach ($colKey in $rowData.Keys) {
        $ws.Range("$colKey$rowKey").Value = [double]$rowData[$colKey]
    }
}